$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 4543.625
$ws.Cells.Item(17, 10).Value = 4543.625
$ws.Cells.Item(17, 12).Value = 13630.875
$ws.Cells.Item(17, 14).Value = -13966.875
$ws.Cells.Item(40, 8).Value = 52237.5
$ws.Cells.Item(40, 9).Value = 36316.668
$ws.Cells.Item(40, 11).Value = 36316.668
$ws.Cells.Item(40, 13).Value = -36141.668
$ws.Cells.Item(53, 8).Value = 651.9091
$ws.Cells.Item(53, 9).Value = 718.2222
$ws.Cells.Item(53, 11).Value = 718.2222
$ws.Cells.Item(53, 13).Value = -81.22220000000004
$ws.Cells.Item(69, 8).Value = 7666.6665
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 13).ClearContents()
$ws.Cells.Item(72, 8).Value = 7666.6665
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 13).ClearContents()
$ws.Cells.Item(93, 8).Value = 53996
$ws.Cells.Item(93, 10).Value = 53996
$ws.Cells.Item(93, 12).Value = 53996
$ws.Cells.Item(93, 14).Value = -58988
$ws.Cells.Item(97, 8).Value = 2524.6667
$ws.Cells.Item(97, 10).Value = 2524.6667
$ws.Cells.Item(97, 12).Value = 7574.000100000001
$ws.Cells.Item(97, 14).Value = -8566.000100000001
$ws.Cells.Item(99, 8).Value = 2971.4614
$ws.Cells.Item(99, 9).Value = 748.44446
$ws.Cells.Item(99, 11).Value = 2245.33338
$ws.Cells.Item(99, 13).Value = -747.33338
$ws.Cells.Item(118, 8).Value = 347.5
$ws.Cells.Item(118, 9).Value = 347.5
$ws.Cells.Item(118, 11).Value = 1042.5
$ws.Cells.Item(118, 13).Value = 614.5
$ws.Cells.Item(125, 8).Value = 2179.5
$ws.Cells.Item(125, 9).Value = 756.4286
$ws.Cells.Item(125, 10).Value = 5500
$ws.Cells.Item(125, 11).Value = 6807.8574
$ws.Cells.Item(125, 12).Value = 49500
$ws.Cells.Item(125, 13).Value = -4347.8574
$ws.Cells.Item(125, 14).Value = -54420
$ws.Cells.Item(132, 8).Value = 10285.56
$ws.Cells.Item(132, 9).Value = 2417.125
$ws.Cells.Item(132, 11).Value = 7251.375
$ws.Cells.Item(132, 13).Value = -4721.375
$ws.Cells.Item(138, 8).Value = 6253.628
$ws.Cells.Item(138, 9).Value = 1425.7142
$ws.Cells.Item(138, 11).Value = 4277.142599999999
$ws.Cells.Item(138, 13).Value = 862.8574000000008
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 21646.256
$ws.Cells.Item(32, 9).Value = 21231.316
$ws.Cells.Item(32, 11).Value = 21231.316
$ws.Cells.Item(32, 13).Value = -20944.316
$ws.Cells.Item(61, 8).Value = 5047.4116
$ws.Cells.Item(61, 9).Value = 4473.773
$ws.Cells.Item(61, 10).Value = 6099.0835
$ws.Cells.Item(61, 11).Value = 4473.773
$ws.Cells.Item(61, 12).Value = 6099.0835
$ws.Cells.Item(61, 13).Value = -4261.773
$ws.Cells.Item(61, 14).Value = -6523.0835
$ws.Cells.Item(132, 8).Value = 24905.633
$ws.Cells.Item(132, 9).Value = 41946.855
$ws.Cells.Item(132, 11).Value = 125840.565
$ws.Cells.Item(132, 13).Value = -123310.565
$ws.Cells.Item(136, 8).Value = 5047.4116
$ws.Cells.Item(136, 9).Value = 4473.773
$ws.Cells.Item(136, 10).Value = 6099.0835
$ws.Cells.Item(136, 11).Value = 13421.319
$ws.Cells.Item(136, 12).Value = 18297.2505
$ws.Cells.Item(136, 13).Value = -10871.319
$ws.Cells.Item(136, 14).Value = -23397.2505
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1043177.8
$ws.Cells.Item(99, 9).Value = 1158808
$ws.Cells.Item(99, 11).Value = 1158808
$ws.Cells.Item(99, 13).Value = -1157310
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 11902.6
$ws.Cells.Item(99, 9).Value = 4012
$ws.Cells.Item(99, 11).Value = 4012
$ws.Cells.Item(99, 13).Value = -2514
$ws.Cells.Item(126, 8).Value = 11902.6
$ws.Cells.Item(126, 9).Value = 4012
$ws.Cells.Item(126, 11).Value = 12036
$ws.Cells.Item(126, 13).Value = -9566
$ws.Cells.Item(134, 8).Value = 2723.6155
$ws.Cells.Item(134, 9).Value = 2373.6843
$ws.Cells.Item(134, 11).Value = 7121.0529
$ws.Cells.Item(134, 13).Value = -4586.0529
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 58571540
$ws.Cells.Item(11, 9).Value = 128.33333
$ws.Cells.Item(11, 11).Value = 384.99999
$ws.Cells.Item(11, 13).Value = -244.99999
$ws.Cells.Item(81, 8).Value = 3874.889
$ws.Cells.Item(81, 9).Value = 1964.25
$ws.Cells.Item(81, 10).Value = 5403.4
$ws.Cells.Item(81, 11).Value = 5892.75
$ws.Cells.Item(81, 12).Value = 16210.2
$ws.Cells.Item(81, 13).Value = -4769.75
$ws.Cells.Item(81, 14).Value = -18456.2
$ws.Cells.Item(84, 8).Value = 3874.889
$ws.Cells.Item(84, 9).Value = 1964.25
$ws.Cells.Item(84, 10).Value = 5403.4
$ws.Cells.Item(84, 11).Value = 17678.25
$ws.Cells.Item(84, 12).Value = 48630.6
$ws.Cells.Item(84, 13).Value = -12062.25
$ws.Cells.Item(84, 14).Value = -59862.6
$ws.Cells.Item(132, 8).Value = 1226.25
$ws.Cells.Item(132, 9).Value = 1012.94116
$ws.Cells.Item(132, 11).Value = 9116.470439999999
$ws.Cells.Item(132, 13).Value = -6586.470439999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(3, 8).Value = 50000000
$ws.Cells.Item(3, 9).Value = 50000000
$ws.Cells.Item(3, 11).Value = 50000000
$ws.Cells.Item(3, 13).Value = -49999884
$ws.Cells.Item(93, 8).Value = 28996
$ws.Cells.Item(93, 10).Value = 28998.25
$ws.Cells.Item(93, 12).Value = 28998.25
$ws.Cells.Item(93, 14).Value = -32742.25
$ws.Cells.Item(97, 8).Value = 683.94116
$ws.Cells.Item(97, 9).Value = 612.61536
$ws.Cells.Item(97, 10).Value = 915.75
$ws.Cells.Item(97, 11).Value = 612.61536
$ws.Cells.Item(97, 12).Value = 915.75
$ws.Cells.Item(97, 13).Value = -116.61536
$ws.Cells.Item(97, 14).Value = -1907.75
$ws.Cells.Item(122, 8).Value = 631593.0600000001
$ws.Cells.Item(122, 9).Value = 1670166.5
$ws.Cells.Item(122, 11).Value = 5010499.5
$ws.Cells.Item(122, 13).Value = -5008049.5
$ws.Cells.Item(132, 8).Value = 4588.3335
$ws.Cells.Item(132, 9).Value = 4123.9473
$ws.Cells.Item(132, 11).Value = 12371.8419
$ws.Cells.Item(132, 13).Value = -9841.841899999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3983.1392
$ws.Cells.Item(7, 9).Value = 3341.8909
$ws.Cells.Item(7, 10).Value = 5452.6665
$ws.Cells.Item(7, 11).Value = 3341.8909
$ws.Cells.Item(7, 12).Value = 5452.6665
$ws.Cells.Item(7, 13).Value = -3229.8909
$ws.Cells.Item(7, 14).Value = -5676.6665
$ws.Cells.Item(126, 8).Value = 3983.1392
$ws.Cells.Item(126, 9).Value = 3341.8909
$ws.Cells.Item(126, 10).Value = 5452.6665
$ws.Cells.Item(126, 11).Value = 10025.6727
$ws.Cells.Item(126, 12).Value = 16357.9995
$ws.Cells.Item(126, 13).Value = -7555.672699999999
$ws.Cells.Item(126, 14).Value = -21297.9995
$ws.Cells.Item(136, 8).Value = 3819.6948
$ws.Cells.Item(136, 9).Value = 2792.8635
$ws.Cells.Item(136, 11).Value = 8378.5905
$ws.Cells.Item(136, 13).Value = -5828.5905
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(102, 8).Value = 30000
$ws.Cells.Item(102, 9).Value = 30000
$ws.Cells.Item(102, 11).Value = 30000
$ws.Cells.Item(102, 13).Value = -26755
$ws.Cells.Item(125, 8).Value = 60000
$ws.Cells.Item(125, 10).Value = 60000
$ws.Cells.Item(125, 12).Value = 60000
$ws.Cells.Item(125, 14).Value = -69840
$ws.Cells.Item(135, 8).Value = 108500
$ws.Cells.Item(135, 10).Value = 108500
$ws.Cells.Item(135, 12).Value = 108500
$ws.Cells.Item(135, 14).Value = -118640
$ws.Cells.Item(136, 8).Value = 3215.55
$ws.Cells.Item(136, 10).Value = 7256.625
$ws.Cells.Item(136, 12).Value = 21769.875
$ws.Cells.Item(136, 14).Value = -26869.875
